$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing existing rows 15-101 down to 16-102.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly price record.
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C15").Value = "Los Lagos"
$ws.Range("D15").Value = 44831
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 100112031
$ws.Range("G15").Value = "Poroto verde"
$ws.Range("H15").Value = "Magnum"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 35
$ws.Range("K15").Value = 33000
$ws.Range("L15").Value = 33000
$ws.Range("M15").Value = 33000
$ws.Range("N15").Value = "`$/malla 25 kilos"
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 1320
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"

# Ensure the date cell carries the same date number format as the rest of column D.
$ws.Range("D15").NumberFormat = $ws.Range("D16").NumberFormat
